# Atualização automática de DOIS_IRMAOS.xlsx
#
# 1. Rename "Paineis DARQ" -> "PAINEIS DARQ"
# 2. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3. Delete sheet "Desarquivamentos Pendentes"
# (DGC sheet content is unaffected - only renumbered indices after the
#  other sheet's removal, which the engine handles automatically.)

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Paineis DARQ")
$ws1.Name = "PAINEIS DARQ"

$ws2 = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$ws2.Name = "RECOLHIMENTO X ELIMINAÇÃO"

$ws3 = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$ws3.Delete()

$excel.DisplayAlerts = $true
